# Sign off the timesheet:
#  - Supervisor Name (G6) on the "Weekly" sheet
#  - Supervisor initials + sign-off date (A27 / D27) in the signature block
#  - Selection left on I27, mirroring where the signer was last working

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor name, formatted like the existing "Employee Name" entry (G4).
$ws.Range("G4").Copy()
$ws.Range("G6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor sign-off row: initials + date, formatted like the employee's
# own sign-off row above it (A25 / D25).
$ws.Range("A25").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = "P.S"

$ws.Range("D25").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 41698

[void]$ws.Range("I27").Select()
